$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.216.56"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.570.42"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.492"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("E7").Value = "  +0.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.09"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.793.03"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.560.22"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.211.63"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.38"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0704"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.43"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "216.38"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.36%  "
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("E22").Value = "  +1.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.25"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.09"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.71"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.11"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.107"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.97%  "
$ws.Range("E30").Value = "  +2.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0474"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.98%  "
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.454.03"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.02%  "
$ws.Range("E34").Value = "  +1.90%  "
$ws.Range("E35").Value = "  +5.84%  "
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("E37").Value = "  +1.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0168"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.82"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.68%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.80"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.707.21"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.93"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.92%  "
$ws.Range("E49").Value = "  +3.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0521"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0962"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.59%  "
